# Update the "scraped_at" timestamps (column K) on the "snapshot" sheet.
# These values represent a fresh scrape re-run, so every existing
# timestamp in K2:K42 is replaced with a newer one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

$ws.Range("K2").Value = "2025-11-20T12:11:09.853953+00:00"
$ws.Range("K3").Value = "2025-11-20T12:11:09.853986+00:00"
$ws.Range("K4").Value = "2025-11-20T12:11:09.854005+00:00"
$ws.Range("K5").Value = "2025-11-20T12:11:11.823315+00:00"
$ws.Range("K6").Value = "2025-11-20T12:11:11.823345+00:00"
$ws.Range("K7").Value = "2025-11-20T12:11:13.686998+00:00"
$ws.Range("K8").Value = "2025-11-20T12:11:16.341718+00:00"
$ws.Range("K9").Value = "2025-11-20T12:11:18.608285+00:00"
$ws.Range("K10").Value = "2025-11-20T12:11:18.608314+00:00"
$ws.Range("K11").Value = "2025-11-20T12:11:18.608334+00:00"
$ws.Range("K12").Value = "2025-11-20T12:11:21.275996+00:00"
$ws.Range("K13").Value = "2025-11-20T12:11:23.996460+00:00"
$ws.Range("K14").Value = "2025-11-20T12:11:26.728213+00:00"
$ws.Range("K15").Value = "2025-11-20T12:11:28.928182+00:00"
$ws.Range("K16").Value = "2025-11-20T12:11:28.928211+00:00"
$ws.Range("K17").Value = "2025-11-20T12:11:28.928228+00:00"
$ws.Range("K18").Value = "2025-11-20T12:11:31.165844+00:00"
$ws.Range("K19").Value = "2025-11-20T12:11:33.471252+00:00"
$ws.Range("K20").Value = "2025-11-20T12:11:36.155866+00:00"
$ws.Range("K21").Value = "2025-11-20T12:11:36.155898+00:00"
$ws.Range("K22").Value = "2025-11-20T12:11:38.390554+00:00"
$ws.Range("K23").Value = "2025-11-20T12:11:38.390583+00:00"
$ws.Range("K24").Value = "2025-11-20T12:11:38.390600+00:00"
$ws.Range("K25").Value = "2025-11-20T12:11:40.978804+00:00"
$ws.Range("K26").Value = "2025-11-20T12:11:40.978833+00:00"
$ws.Range("K27").Value = "2025-11-20T12:11:43.601049+00:00"
$ws.Range("K28").Value = "2025-11-20T12:11:43.601080+00:00"
$ws.Range("K29").Value = "2025-11-20T12:11:43.601101+00:00"
$ws.Range("K30").Value = "2025-11-20T12:11:45.919960+00:00"
$ws.Range("K31").Value = "2025-11-20T12:11:45.919992+00:00"
$ws.Range("K32").Value = "2025-11-20T12:11:48.210729+00:00"
$ws.Range("K33").Value = "2025-11-20T12:11:48.210761+00:00"
$ws.Range("K34").Value = "2025-11-20T12:11:48.210782+00:00"
$ws.Range("K35").Value = "2025-11-20T12:11:48.210803+00:00"
$ws.Range("K36").Value = "2025-11-20T12:11:48.210819+00:00"
$ws.Range("K37").Value = "2025-11-20T12:11:50.520050+00:00"
$ws.Range("K38").Value = "2025-11-20T12:11:50.520085+00:00"
$ws.Range("K39").Value = "2025-11-20T12:11:55.580757+00:00"
$ws.Range("K40").Value = "2025-11-20T12:11:55.580787+00:00"
$ws.Range("K41").Value = "2025-11-20T12:11:55.580805+00:00"
$ws.Range("K42").Value = "2025-11-20T12:11:57.822736+00:00"
